# Update COVID-19 "paises" data: refresh case counts for several countries
# and re-sort (by total cases, column B, descending) the groups of rows
# whose order changed as a result, plus bump the "last updated" timestamp.
#
# Column layout (row 3 header):
#   A: Pais            B: Casos totales   C: Nuevos casos   D: Casos activos
#   E: Recuperados     F: Casos criticos  G: Muertes hoy    H: Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp in the footer row ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 00:20"

# --- Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 161044
$ws.Cells.Item(4, 3).Value = 17553
$ws.Cells.Item(4, 4).Value = 5245
$ws.Cells.Item(4, 5).Value = 152831
$ws.Cells.Item(4, 7).Value = 385
$ws.Cells.Item(4, 8).Value = 2968

# --- Alemania (row 8) ---
$ws.Cells.Item(8, 2).Value = 66885
$ws.Cells.Item(8, 3).Value = 4450
$ws.Cells.Item(8, 5).Value = 52740

# --- Canada (row 18) ---
$ws.Cells.Item(18, 5).Value = 6248
$ws.Cells.Item(18, 7).Value = 21
$ws.Cells.Item(18, 8).Value = 86

# --- Peru (row 48) ---
$ws.Cells.Item(48, 4).Value = 53
$ws.Cells.Item(48, 5).Value = 873
$ws.Cells.Item(48, 7).Value = 6
$ws.Cells.Item(48, 8).Value = 24

# --- Rows 73-76: Bulgaria / Bosnia y Herzegovina / Eslovaquia / Principado de
#     Andorra get updated totals and are re-sorted (Andorra & Bosnia jump
#     ahead of Bulgaria & Eslovaquia, whose own figures are unchanged) ---
$ws.Cells.Item(73, 1).Value = "Principado de Andorra"
$ws.Cells.Item(73, 2).Value = 370
$ws.Cells.Item(73, 3).Value = 36
$ws.Cells.Item(73, 4).Value = 10
$ws.Cells.Item(73, 5).Value = 352
$ws.Cells.Item(73, 6).Value = 10
$ws.Cells.Item(73, 7).Value = 2

$ws.Cells.Item(74, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(74, 2).Value = 368
$ws.Cells.Item(74, 3).Value = 45
$ws.Cells.Item(74, 5).Value = 341
$ws.Cells.Item(74, 7).Value = 4
$ws.Cells.Item(74, 8).Value = 10

$ws.Cells.Item(75, 1).Value = "Bulgaria"
$ws.Cells.Item(75, 2).Value = 359
$ws.Cells.Item(75, 3).Value = 13
$ws.Cells.Item(75, 4).Value = 17
$ws.Cells.Item(75, 5).Value = 334
$ws.Cells.Item(75, 6).Value = 13
$ws.Cells.Item(75, 8).Value = 8

$ws.Cells.Item(76, 1).Value = "Eslovaquia"
$ws.Cells.Item(76, 2).Value = 336
$ws.Cells.Item(76, 3).Value = 22
$ws.Cells.Item(76, 4).Value = 7
$ws.Cells.Item(76, 5).Value = 329
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 8).Value = 0

# --- Rows 96-97: Costa de Marfil overtakes Islas Feroe ---
$ws.Cells.Item(96, 1).Value = "Costa de Marfil"
$ws.Cells.Item(96, 3).Value = 3
$ws.Cells.Item(96, 4).Value = 6
$ws.Cells.Item(96, 5).Value = 161
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 8).Value = 1

$ws.Cells.Item(97, 1).Value = "Islas Feroe"
$ws.Cells.Item(97, 2).Value = 168
$ws.Cells.Item(97, 3).Value = 9
$ws.Cells.Item(97, 4).Value = 70
$ws.Cells.Item(97, 5).Value = 98
$ws.Cells.Item(97, 6).Value = 1
$ws.Cells.Item(97, 8).Value = 0

# --- Rows 134-137: Polinesia Francesa overtakes Jamaica/Guatemala/Zambia ---
$ws.Cells.Item(134, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(134, 3).Value = 6
$ws.Cells.Item(134, 4).Value = 0
$ws.Cells.Item(134, 5).Value = 36
$ws.Cells.Item(134, 6).Value = 2
$ws.Cells.Item(134, 8).Value = 0

$ws.Cells.Item(135, 1).Value = "Jamaica"
$ws.Cells.Item(135, 3).Value = 4
$ws.Cells.Item(135, 4).Value = 2
$ws.Cells.Item(135, 5).Value = 33
$ws.Cells.Item(135, 6).Value = 0

$ws.Cells.Item(136, 1).Value = "Guatemala"
$ws.Cells.Item(136, 2).Value = 36
$ws.Cells.Item(136, 3).Value = 2
$ws.Cells.Item(136, 4).Value = 10
$ws.Cells.Item(136, 5).Value = 25
$ws.Cells.Item(136, 6).Value = 1
$ws.Cells.Item(136, 8).Value = 1

$ws.Cells.Item(137, 1).Value = "Zambia"
$ws.Cells.Item(137, 3).Value = 6
$ws.Cells.Item(137, 6).Value = 0

# --- Rows 144-147: Niger overtakes Mali/Etiopia/Guinea ---
$ws.Cells.Item(144, 1).Value = "Niger"
$ws.Cells.Item(144, 2).Value = 27
$ws.Cells.Item(144, 3).Value = 9
$ws.Cells.Item(144, 5).Value = 24
$ws.Cells.Item(144, 7).Value = 2
$ws.Cells.Item(144, 8).Value = 3

$ws.Cells.Item(145, 1).Value = "Mali"
$ws.Cells.Item(145, 2).Value = 25
$ws.Cells.Item(145, 3).Value = 7
$ws.Cells.Item(145, 4).Value = 0
$ws.Cells.Item(145, 5).Value = 23
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 2

$ws.Cells.Item(146, 1).Value = "Etiopia"
$ws.Cells.Item(146, 2).Value = 23
$ws.Cells.Item(146, 3).Value = 2
$ws.Cells.Item(146, 4).Value = 4
$ws.Cells.Item(146, 5).Value = 19
$ws.Cells.Item(146, 6).Value = 1

$ws.Cells.Item(147, 1).Value = "Guinea"
$ws.Cells.Item(147, 3).Value = 6
$ws.Cells.Item(147, 5).Value = 22
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 0

# --- Rows 165-169: Siria overtakes Groenlandia/Suazilandia/Granada/Santa Lucia ---
$ws.Cells.Item(165, 1).Value = "Siria"
$ws.Cells.Item(165, 3).Value = 1
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 7).Value = 1
$ws.Cells.Item(165, 8).Value = 2

$ws.Cells.Item(166, 1).Value = "Groenlandia"
$ws.Cells.Item(166, 2).Value = 10
$ws.Cells.Item(166, 4).Value = 2
$ws.Cells.Item(166, 5).Value = 8

$ws.Cells.Item(167, 1).Value = "Suazilandia"

$ws.Cells.Item(168, 1).Value = "Granada"
$ws.Cells.Item(168, 4).Value = 0
$ws.Cells.Item(168, 5).Value = 9

$ws.Cells.Item(169, 1).Value = "Santa Lucia"
$ws.Cells.Item(169, 4).Value = 1
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 7).Value = 0
$ws.Cells.Item(169, 8).Value = 0
